# Updated cryptos list data: apply per-cell changes as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.044.01'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.280.38'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''317.43'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = '''100.83'
$ws.Range('E6').Value = '  -4.37%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '''0.600'
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').Value = '''38.93'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').Value = '''0.0898'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = '''8.20'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '''0.951'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = '''15.09'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').Value = '2.626.59'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '2.287.56'
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('D18').Value = '42.058.13'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '''7.33'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').Value = '''12.84'
$ws.Range('E21').Value = '  +31.72%  '
$ws.Range('D22').Value = '''72.51'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '''3.52'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '''266.44'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('E25').Value = '  -5.03%  '
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').Value = '''10.75'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '''2.32'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('D29').Value = '''22.36'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').Value = '''37.15'
$ws.Range('E30').Value = '  +2.75%  '
$ws.Range('D31').Value = '''165.07'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '''6.05'
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('D33').Value = '''0.0865'
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('D34').Value = '''0.132'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').Value = '''2.56'
$ws.Range('E35').Value = '  -12.38%  '
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('D37').Value = '''4.55'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').Value = '''0.0353'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('E39').Value = '  -5.46%  '
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').Value = '''1.51'
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = '''67.95'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').Value = '''91.16'
$ws.Range('E45').Value = '  -7.88%  '
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('D47').Value = '''11.81'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').Value = '''78.29'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = '''8.90'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.598.06'
$ws.Range('E50').Value = '  +3.39%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''5.19'
$ws.Range('E51').Value = '  -2.25%  '
